# Still working on the trafo model.
# Applies:
#   - busbars: add Bus3 (row 4)
#   - lines:   add Line1 (row 2)
#   - loads:   update "Load 1" row values (p_nom_mw/q_nom_mvar/bus_idx + text "1.0" flag)
#   - trafos:  add tap-changer columns (tap_pos, tap_change, tap_min, tap_max) + new values

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "busbars": new Bus3 row
# ---------------------------------------------------------------------------
$wsBus = $wb.Worksheets.Item("busbars")
$wsBus.Range("A4").Value = 2
$wsBus.Range("B4").Value = "Bus3"
$wsBus.Range("C4").Value = 22
$wsBus.Range("D4").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "lines": new Line1 row
# ---------------------------------------------------------------------------
$wsLines = $wb.Worksheets.Item("lines")
$wsLines.Range("A2").Value = "Line1"
$wsLines.Range("B2").Value = 22
$wsLines.Range("C2").Value = 1
$wsLines.Range("D2").Value = 1
$wsLines.Range("E2").Value = 4
$wsLines.Range("F2").Value = 0
$wsLines.Range("G2").Value = 1
$wsLines.Range("H2").Value = 2
$wsLines.Range("I2").Value = 0
$wsLines.Range("A3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "loads": update the "Load 1" row
# ---------------------------------------------------------------------------
$wsLoads = $wb.Worksheets.Item("loads")
# B2 must become the literal text "1.0" (not the number 1) while keeping the
# default "Normal" cell style (no numberformat residue left behind).
$wsLoads.Range("B2").Value = "'1.0"
$wsLoads.Range("B2").Style = "Normal"
$wsLoads.Range("C2").Value = 10
$wsLoads.Range("D2").Value = 5
$wsLoads.Range("E2").Value = 2
$wsLoads.Range("A3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "trafos": add tap-changer columns + values for T1
# ---------------------------------------------------------------------------
$wsTrafos = $wb.Worksheets.Item("trafos")
$wsTrafos.Range("K1").Value = "tap_pos"
$wsTrafos.Range("L1").Value = "tap_change"
$wsTrafos.Range("M1").Value = "tap_min"
$wsTrafos.Range("N1").Value = "tap_max"

$wsTrafos.Range("E2").Value = 0.03
$wsTrafos.Range("H2").Value = 0.001
$wsTrafos.Range("K2").Value = 1
$wsTrafos.Range("L2").Value = 0.01
$wsTrafos.Range("M2").Value = -5
$wsTrafos.Range("N2").Value = 5

# Keep "trafos" as the active/selected tab, matching the saved workbook view.
$wsTrafos.Range("K4").Select() | Out-Null
